$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 headers to reflect new column layout (Watch and Sec. Position inserted,
# Potential and Profile removed) starting at column J.
$ws.Range("J1").Value = "Watch"
$ws.Range("K1").Value = "Height"
$ws.Range("L1").Value = "End Contract"
$ws.Range("M1").Value = "Market Value"
$ws.Range("N1").Value = "Position"
$ws.Range("O1").Value = "Sec. Position"
$ws.Range("P1").Value = "Nationality"
$ws.Range("Q1").Value = "Agent"
$ws.Range("R1").Value = "Performance"

# Update the active selection to match the saved view state.
$ws.Range("M13").Select()
